$wb = $excel.ActiveWorkbook

# 1. Rename "Foglio1" -> "rapporti richieste posti"
$ws3 = $wb.Worksheets.Item("Foglio1")
$ws3.Name = "rapporti richieste posti"

# 2. Add a new sheet after it, named "Foglio2"
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Foglio2"

# 3. Populate Foglio2
$ws4.Range("A1").Value = "uni"
$ws4.Range("B1").Value = "corsi possibili"
$ws4.Range("A2").Value = "UNIVERSITAT BIELEFELD"

